# Update the "想去人数" (F column) figures on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Map of row -> new F value, per worksheet.
$updates1 = @{
    2  = 50
    3  = 778
    4  = 37
    6  = 60
    8  = 3834
    9  = 85
    10 = 4517
    11 = 492
    12 = 1136
    13 = 65
}

$updates4 = @{
    2  = 50
    3  = 778
    4  = 37
    6  = 60
    9  = 3834
    10 = 85
    11 = 4517
    12 = 492
    13 = 1136
    14 = 65
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates1.Keys) {
    $ws1.Range("F$row").Value = $updates1[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
